$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 328-329, shifting the existing rows 328-408 down to 330-410.
$ws.Rows("328:329").Insert()

# Populate the first new row (328)
$ws.Range("A328").Value2 = 10
$ws.Range("B328").Value2 = "Vega Modelo de Temuco"
$ws.Range("C328").Value2 = "La Araucanía"
$ws.Range("D328").Value2 = 44785
$ws.Range("E328").Value2 = 9
$ws.Range("F328").Value2 = 100112040
$ws.Range("G328").Value2 = "Cilantro"
$ws.Range("H328").Value2 = "Sin especificar"
$ws.Range("I328").Value2 = "Primera"
$ws.Range("J328").Value2 = 30
$ws.Range("K328").Value2 = 5000
$ws.Range("L328").Value2 = 5000
$ws.Range("M328").Value2 = 5000
$ws.Range("N328").Value2 = "`$/docena de atados (2 kilos)"
$ws.Range("O328").Value2 = "Provincia de Cautín"
$ws.Range("P328").Value2 = 2500
$ws.Range("Q328").Value2 = 2
$ws.Range("R328").Value2 = "Hortaliza"

# Populate the second new row (329)
$ws.Range("A329").Value2 = 10
$ws.Range("B329").Value2 = "Vega Modelo de Temuco"
$ws.Range("C329").Value2 = "La Araucanía"
$ws.Range("D329").Value2 = 44785
$ws.Range("E329").Value2 = 9
$ws.Range("F329").Value2 = 100112040
$ws.Range("G329").Value2 = "Cilantro"
$ws.Range("H329").Value2 = "Sin especificar"
$ws.Range("I329").Value2 = "Primera"
$ws.Range("J329").Value2 = 40
$ws.Range("K329").Value2 = 3500
$ws.Range("L329").Value2 = 3500
$ws.Range("M329").Value2 = 3500
$ws.Range("N329").Value2 = "`$/docena de atados (2 kilos)"
$ws.Range("O329").Value2 = "Región Metropolitana"
$ws.Range("P329").Value2 = 1750
$ws.Range("Q329").Value2 = 2
$ws.Range("R329").Value2 = "Hortaliza"
